$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the stray "_GoBack" bookmark that used to sit right after
# "Name: Nathan Cochrane" in the first paragraph.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# Change 2: split the "Weak points..." run right after the apostrophe in
# "I'" and drop a new "_GoBack" bookmark (collapsed, zero-length) at that
# split point - this is where the author's cursor was left last.
# ---------------------------------------------------------------------------
$weak = $d.Content
[void]$weak.Find.Execute("Weak points and what I" + [char]8217, $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$splitPoint = $weak.End
$markRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $markRange)

# ---------------------------------------------------------------------------
# Change 3: append a new sentence after "... to develop." in the
# "Challenges in the context of the project:" paragraph.  Each chunk is
# inserted at an advancing, collapsed point so Word starts a fresh run for
# each piece (mirroring the separate <w:r> runs introduced upstream).
# ---------------------------------------------------------------------------
$tail = $d.Content
[void]$tail.Find.Execute("to develop.", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$tail.Collapse(0)  # wdCollapseEnd

$tail.InsertAfter(" This means that we" + [char]8217 + "ll need to be very ")
$tail.Collapse(0)

$wordChunk = $d.Range($tail.End, $tail.End)
$wordChunk.InsertAfter("organised")

$period = $d.Range($wordChunk.End, $wordChunk.End)
$period.InsertAfter(".")
